# fe_table_poster.docx: update/clear the "TWFE + Covariates" column values.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2: Log real minimum wage -> coefficient 1.151 becomes 1.205
$t.Cell(2, 2).Range.Text = "1.205"

# Row 3: corresponding standard error (1.538) becomes (1.632)
$t.Cell(3, 2).Range.Text = "(1.632)"

# The remaining covariate coefficient/SE pairs in this column are cleared
# out entirely (text removed and the paragraph alignment reset to the
# style default, which drops the explicit <w:jc> left-alignment override).
$rowsToClear = @(4, 5, 6, 7, 8, 9, 10, 11)
foreach ($r in $rowsToClear) {
    $cell = $t.Cell($r, 2)
    $cell.Range.Delete()
    $cell.Range.ParagraphFormat.Alignment = 0
}
